$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.393.09"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.878.08"
$ws.Range("D3").ClearFormats()

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7161"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.70"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07979"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3145"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.90"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08081"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.873.20"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.67"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.226"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7075"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.395"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008443"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.398.44"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.85"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.34"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.134.64"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.670"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.65%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1576"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.063"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.90"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.98"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.509"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.416"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.314"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.224"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05308"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.942"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7577"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.174"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.702"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("E38").Value = "  +0.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.275.35"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.756"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.406"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9072"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "74.11"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "111.56"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.65%  "

$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.029.90"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.806"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5209"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.521"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4340"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.37%  "
